# fix bugs result double display
# - Row 34: a new leading "jenis/unit" column is inserted (values shift one column right).
# - Rows 37-64 & 66: a duplicate "D" amount column is added next to the existing "C" column
#   (mirrors the value) and both get an integer ("0") number format.
# - Rows 70/71/74/75 (two regex-lookup rows + two raw numbers that were rendering as a
#   doubled-up display) are replaced by a clean two-column lookup table of crop names
#   (agroforestry variants in B73:B89, plain crop names in D73:D93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------------
# Row 34 : shift all existing values one column to the right (A->B, B->C, C->D, ...)
#          and put the new "utama" marker in the now-empty A34.
# ----------------------------------------------------------------------------------
$row34 = @(
    "utama",
    "tbs",
    "kg",
    0,
    0,
    0,
    2658.5160642857145,
    4044.8571428571431,
    4829.469542857144,
    6852.8980928571427,
    8728.8017142857152,
    10457.180407142858,
    12038.034171428573,
    13471.363007142858,
    14757.166914285717,
    15895.445892857144,
    16886.19994285714,
    17729.42906428572,
    18425.133257142861,
    18973.312521428572,
    19373.966857142863,
    19627.096264285716,
    19732.700742857145,
    19690.780292857147,
    19501.334914285715,
    19164.364607142859,
    18679.869371428573,
    18047.849207142859
)
$arr34 = New-Object 'object[,]' 1,$row34.Length
for ($i = 0; $i -lt $row34.Length; $i++) {
    $arr34[0,$i] = $row34[$i]
}
$ws.Range("A34:AB34").Value = $arr34
$ws.Range("A34").VerticalAlignment = -4108

# ----------------------------------------------------------------------------------
# Rows 37-64 : duplicate column C into a new column D (same value), both formatted as
#              integers ("0" number format, cellXfs style index 10 in the saved file).
# ----------------------------------------------------------------------------------
$cvals = [ordered]@{
    37 = 9000
    38 = 14000
    39 = 1000
    40 = 1000
    41 = 90000
    42 = 102000
    43 = 95000
    44 = 100000
    45 = 38000
    46 = 93000
    47 = 585000
    48 = 500000
    49 = 460000
    50 = 80000
    51 = 5500000
    52 = 50000
    53 = 300000
    54 = 200000
    55 = 100000
    56 = 30000
    57 = 100000
    58 = 100000
    59 = 100000
    60 = 100000
    61 = 100000
    62 = 100000
    63 = 100000
    64 = 100000
}
foreach ($r in $cvals.Keys) {
    $v = $cvals[$r]
    $ws.Cells.Item($r, 3).Value = $v
    $ws.Cells.Item($r, 4).Value = $v
}
$ws.Range("C37:D64").NumberFormat = "0"

# Row 66 : same duplication, but no special number format on this pair.
$ws.Cells.Item(66, 4).Value = 1570

# ----------------------------------------------------------------------------------
# Rows 70/71/74/75 used to hold two header/regex rows and two bare numbers that
# ended up displayed twice. Clear them and replace with the real lookup table.
# ----------------------------------------------------------------------------------
$ws.Cells.Item(70, 2).ClearContents()
$ws.Cells.Item(71, 2).ClearContents()
$ws.Cells.Item(74, 2).ClearContents()
$ws.Cells.Item(75, 2).ClearContents()

# Column B (rows 73-89): agroforestry crop variants.
$bVals = @(
    "COKLAT AF",
    "DUKU AF",
    "JATI AF",
    "KARET AF",
    "KARET AF PADI BUAH",
    "KARET AF PADI KOPI",
    "KAYU MANIS AF",
    "KEBUN CAMPUR AF",
    "KELAPA AF",
    "KELAPA AF COKLAT",
    "KELAPA AF JAGUNG",
    "KELAPA SAWIT AF",
    "KEMIRI AF",
    "KOPI AF",
    "LOGGING AF",
    "SALAK AF",
    "SENGON AF"
)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws.Cells.Item(73 + $i, 2).Value = $bVals[$i]
}

# Column D (rows 73-93): plain crop names.
$dVals = @(
    "CENGKEH",
    "COKLAT",
    "JAGUNG",
    "JATI",
    "KARET",
    "KAYU MANIS  ",
    "KELAPA",
    "KELAPA ",
    "KELAPA SAWIT",
    "KELAPA SAWIT LARGE SCALE",
    "KENTANG",
    "KOPI   ",
    "LADA",
    "PADI",
    "PADI DRYLAND",
    "PADI IRIGASI",
    "PADI WETLAND",
    "PINANG",
    "SAGU",
    "SENGON",
    "TEH"
)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $ws.Cells.Item(73 + $i, 4).Value = $dVals[$i]
}

# ----------------------------------------------------------------------------------
# Restore the on-screen selection to the new table's first row.
# ----------------------------------------------------------------------------------
$ws.Range("B66:D66").Select()
